# labor-timekeeper simulator: full-month coverage, persist logs, fix employees
$wb = $excel.ActiveWorkbook

$tsSheet = $wb.Worksheets.Item("Weekly Timesheet")
$jsSheet = $wb.Worksheets.Item("Jason Schema")

# --- fix employee/client name: "Richer" -> "Evans" (shared string referenced
# from both sheets; update every referencing cell so the old string drops out) ---
$tsSheet.Range("B5").Value = "Evans"
$jsSheet.Range("D5").Value = "Evans"

# --- fix employee id ---
$jsSheet.Range("B2:B5").Value = "emp_jp4mlvog"

# --- Weekly Timesheet: populate Rate / Total for each logged day (full-month
# simulator coverage now persists the $140/hr rate and computed totals
# instead of leaving them at 0) ---
$tsSheet.Range("E2").Value = 140
$tsSheet.Range("F2").Value = 980

$tsSheet.Range("E3").Value = 140
$tsSheet.Range("F3").Value = 1120

$tsSheet.Range("E4").Value = 140
$tsSheet.Range("F4").Value = 2800

$tsSheet.Range("E5").Value = 140
$tsSheet.Range("F5").Value = 2800

# subtotal / admin subtotal / grand total rollups
$tsSheet.Range("F7").Value = 7700
$tsSheet.Range("F11").Value = 7700
$tsSheet.Range("F12").Value = 7700

# --- Jason Schema: mirror the same Rate / Total figures ---
$jsSheet.Range("F2").Value = 140
$jsSheet.Range("G2").Value = 980

$jsSheet.Range("F3").Value = 140
$jsSheet.Range("G3").Value = 1120

$jsSheet.Range("F4").Value = 140
$jsSheet.Range("G4").Value = 2800

$jsSheet.Range("F5").Value = 140
$jsSheet.Range("G5").Value = 2800
